$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(5)
$tr = $shape.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.Text = "Our project is a console application, that make it easy and understandable for the people who are learning new information about the basic Finance Challenges, so they can have better general knowledge."
